$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.06379246711731
$ws.Range("B1").Value = 2.957987308502197
$ws.Range("C1").Value = 4.634864807128906
$ws.Range("D1").Value = 1.016380786895752
$ws.Range("E1").Value = 1.234916090965271
